# Weekly Report 13 - update for the week of the 9th October.
#
# 1. Remove the stray "_GoBack" bookmark that was sitting in the title
#    paragraph (an artifact Word leaves at the last edit location).
# 2. Replace the long "We had a meeting on the 29th November..." narrative
#    with the new "We did not have a meeting out of class..." paragraph,
#    collapsing the old multi-run text into a single run.
# 3. Re-create the "_GoBack" bookmark at the new last-edit location, i.e.
#    immediately after the freshly typed text.

$d = $word.ActiveDocument

# --- Step 1: drop the old bookmark from the title line -------------------
$d.Bookmarks.Item("_GoBack").Delete()

# --- Step 2: swap the body text ------------------------------------------
$oldText = "We had a meeting on the 29th November to discuss our plan for the project overall. This meeting was attended by John, Michael and Zach. We concluded that by the 4th October, we should all have our interfaces for the project completed, with suitable pseudocode. On the 6th October, Zach will presented the final version of the rock for the group. On the 11th October, Shane will have completed the rock, enemy and egg models for the project. On the 20th October, Everyone in the group will present a test version of their project. Zach will present a test version of a top down camera that follows the player. John will present a test scenario which shows a working pop-up score. Joseph will present a test of the enemy AI. Shane will present a test version of the egg. Michael will present a test scenario, which will spawn a random number of different blocks that will be stored in a list. Everyone will be finished coding their individual portions of the project by the middle of November. John will then test everybody" + [char]0x2019 + "s code and we will correct any errors in our code."

$newText = "We did not have a meeting out of class. Michael, Rob and Laura did work on getting the game manager to spawn the pickup item. Joseph and Shane worked on getting the AI to move automatically. John and Michael worked on getting the pop up score to work in a test. Zach discussed the code for the camera and rock."

$range = $d.Content
$found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) {
    throw "Could not locate the original meeting paragraph text to replace."
}

# --- Step 3: move "_GoBack" to sit right after the new text --------------
# The engine mis-resolves a zero-width range collapsed exactly at a
# paragraph's end boundary, so nudge past it with a temporary marker
# character, add the bookmark while the position is mid-paragraph, then
# remove the marker again.
$endRange = $d.Content
$endRange.Find.Execute($newText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endRange.Collapse(0)
$endRange.InsertAfter("~")

$markRange = $d.Content
$markRange.Find.Execute($newText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $markRange)

$markerRange = $d.Content
$markerRange.Find.Execute("~", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markerRange.Text = ""
